$d = $word.ActiveDocument

# Insert new sentence about removing Sub-Classification 'Other' rows
$ok1 = $d.Content.Find.Execute("This decreased the total data entries to 197,229. After this", $true, $false, $false, $false, $false, $true, 1, $false, "This decreased the total data entries to 197,229. In addition, rows with a Sub-Classification value of ‘Other’ was removed. This is because ‘Other’ does not tell us anything useful about the data. It decreased the total data entries to 184,841. After this", 2)
if (-not $ok1) { Write-Host "WARNING: replace 1 (Insert new sentence about removing Sub-Classification 'Other' rows) failed" }

# Duplicate entries found: 1,190 -> 1,039
$ok2 = $d.Content.Find.Execute(" revealed 1,190 duplicate entries. ", $true, $false, $false, $false, $false, $true, 1, $false, " revealed 1,039 duplicate entries. ", 2)
if (-not $ok2) { Write-Host "WARNING: replace 2 (Duplicate entries found: 1,190 -> 1,039) failed" }

# Total entries after duplicate removal: 196,039 -> 183,802
$ok3 = $d.Content.Find.Execute("196,039", $true, $false, $false, $false, $false, $true, 1, $false, "183,802", 2)
if (-not $ok3) { Write-Host "WARNING: replace 3 (Total entries after duplicate removal: 196,039 -> 183,802) failed" }

# LowestSalary 0k percentage: 19.72% -> 19.53%
$ok4 = $d.Content.Find.Execute("and so nothing was dropped. Additionally, it was found that 19.72% of jobs had a ", $true, $false, $false, $false, $false, $true, 1, $false, "and so nothing was dropped. Additionally, it was found that 19.53% of jobs had a ", 2)
if (-not $ok4) { Write-Host "WARNING: replace 4 (LowestSalary 0k percentage: 19.72% -> 19.53%) failed" }

# HighestSalary 999k percentage: 3.57% -> 4.35%
$ok5 = $d.Content.Find.Execute(" value of 0k, and 3.57% had a ", $true, $false, $false, $false, $false, $true, 1, $false, " value of 0k, and 4.35% had a ", 2)
if (-not $ok5) { Write-Host "WARNING: replace 5 (HighestSalary 999k percentage: 3.57% -> 4.35%) failed" }

# Duplicate entries removed count: 1,190 -> 1039
$ok6 = $d.Content.Find.Execute("columns, 1,190 duplicate entries", $true, $false, $false, $false, $false, $true, 1, $false, "columns, 1039 duplicate entries", 2)
if (-not $ok6) { Write-Host "WARNING: replace 6 (Duplicate entries removed count: 1,190 -> 1039) failed" }

# Rows affected by Area removal from Title: 8,353 -> 8,011
$ok7 = $d.Content.Find.Execute("8,353", $true, $false, $false, $false, $false, $true, 1, $false, "8,011", 2)
if (-not $ok7) { Write-Host "WARNING: replace 7 (Rows affected by Area removal from Title: 8,353 -> 8,011) failed" }

